# Fill in results for contests 38 through 46 (worksheet rows 50-58).
# Columns E,H,K,N,Q,T,W hold the raw scores for each of the 7 players;
# columns D,G,J,M,P,S,V already contain VLOOKUP/RANK formulas that derive
# the scoring-table value for the corresponding raw score column.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$results = @{
    50 = @{ E = 70;  H = 80;  K = 60;  N = 100; Q = 50;  T = 40;  W = 0   }
    51 = @{ E = 0;   H = 80;  K = 50;  N = 60;  Q = 70;  T = 100; W = 40  }
    52 = @{ E = 0;   H = 100; K = 60;  N = 80;  Q = 70;  T = 40;  W = 50  }
    53 = @{ E = 100; H = 50;  K = 60;  N = 0;   Q = 70;  T = 40;  W = 80  }
    54 = @{ E = 50;  H = 100; K = 50;  N = 60;  Q = 70;  T = 0;   W = 80  }
    55 = @{ E = 60;  H = 40;  K = 50;  N = 0;   Q = 70;  T = 100; W = 80  }
    56 = @{ E = 40;  H = 60;  K = 100; N = 70;  Q = 80;  T = 0;   W = 50  }
    57 = @{ E = 40;  H = 0;   K = 60;  N = 70;  Q = 100; T = 80;  W = 50  }
    58 = @{ E = 80;  H = 70;  K = 50;  N = 0;   Q = 100; T = 60;  W = 40  }
}

foreach ($row in $results.Keys) {
    $scores = $results[$row]
    foreach ($col in $scores.Keys) {
        $ws.Range("$col$row").Value = $scores[$col]
    }
}

# Row 54 has a tie between columns E and K (both scored 50), so the
# RANK/VLOOKUP formula can't resolve a single rank; the author replaced
# those two derived cells with the manually-averaged result.
$ws.Range("D54").Value = -17.5
$ws.Range("J54").Value = -17.5

$excel.Calculate()
